# InstallTracker.xlsx update
#
# The "Daily" install counts (column C) for 8/30-8/31/19-ish entries were
# revised, and three more days of data (rows 87-89) were appended to the
# tracker. Columns D (day-over-day delta) and E (7-day rolling average)
# are formula-driven off column C, so re-entering C re-derives D/E for the
# rows that already had formulas (84:86) and we add matching formulas for
# the newly-populated rows (87:89). Columns F/G (poly-2 / poly-3 estimate
# curves) only depend on column B (day index), so they are untouched and
# simply recalc in place.
#
# Finally the sheet's scroll position / active cell are updated to reflect
# where the editor was working when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Revise existing daily counts (rows 84:86) --------------------------
# D68:D86 and E74:E86 are already formula-filled for these rows, so just
# writing the new "Actual" value is enough for them to recompute.
$ws.Range("C84").Value = 1495
$ws.Range("C85").Value = 1541
$ws.Range("C86").Value = 1580

# ---- Append three new days (rows 87:89) ----------------------------------
# Row 87
$ws.Range("C87").Value = 1614
$ws.Range("C87").NumberFormat = $ws.Range("C86").NumberFormat
$ws.Range("D87").Formula = "=C87-C86"
$ws.Range("D87").NumberFormat = $ws.Range("D86").NumberFormat
$ws.Range("E87").Formula = "=(C87-C80)/7"
$ws.Range("E87").NumberFormat = $ws.Range("E86").NumberFormat

# Row 88
$ws.Range("C88").Value = 1643
$ws.Range("C88").NumberFormat = $ws.Range("C86").NumberFormat
$ws.Range("D88").Formula = "=C88-C87"
$ws.Range("D88").NumberFormat = $ws.Range("D86").NumberFormat
$ws.Range("E88").Formula = "=(C88-C81)/7"
$ws.Range("E88").NumberFormat = $ws.Range("E86").NumberFormat

# Row 89
$ws.Range("C89").Value = 1669
$ws.Range("C89").NumberFormat = $ws.Range("C86").NumberFormat
$ws.Range("D89").Formula = "=C89-C88"
$ws.Range("D89").NumberFormat = $ws.Range("D86").NumberFormat
$ws.Range("E89").Formula = "=(C89-C82)/7"
$ws.Range("E89").NumberFormat = $ws.Range("E86").NumberFormat

# ---- Nudge the charts to pick up the new points --------------------------
# The scatter charts already reference Sheet1!$B$2:$B$115 / $C$2:$C$115 /
# $D$2:$D$115 / $E$2:$E$115, i.e. ranges that already cover the new rows,
# so a source refresh just re-reads the now-current worksheet values.
foreach ($co in $ws.ChartObjects()) {
    try {
        $co.Chart.Refresh()
    } catch {
    }
}

# ---- Restore the editor's scroll position / selection --------------------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 2
$aw.ScrollColumn = 1
$null = $ws.Range("H17").Select()
